$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells that need to hold text (not auto-converted numbers/percentages)
$textCells = @(
    "B8","C8","D8","E8",
    "B9","C9","D9","E9",
    "D2","E2","D3","E3","D4","E4","E5","D6","E6","D7","E7",
    "D10","E10","D11","E11","D12","E12","D13","E13","D14","E14","E15",
    "D16","E16","E17","D18","E18","D19","E19","D20","E20","D21","E21",
    "D22","E22","D23","E23","D24","E24","D25","E25","D26","E26",
    "D38","E38","D39","E39","D40","E40","D41","E41","D42","E42",
    "D43","E43","D44","E44","D45","E45","E46","E47","D48","E48",
    "D49","E49","D50","E50"
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2
$ws.Range("D2").Value = "327.38"
$ws.Range("E2").Value = "0.75%"

# Row 3
$ws.Range("D3").Value = "43.98"
$ws.Range("E3").Value = "-0.34%"

# Row 4
$ws.Range("D4").Value = "5.503"
$ws.Range("E4").Value = "-0.15%"

# Row 5
$ws.Range("E5").Value = "-0.18%"

# Row 6
$ws.Range("D6").Value = "1.999"
$ws.Range("E6").Value = "5.45%"

# Row 7
$ws.Range("D7").Value = "4.319"
$ws.Range("E7").Value = "-0.46%"

# Row 8 (swapped with row 9 in the source data)
$ws.Range("B8").Value = "BTSEToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D8").Value = "2.570"
$ws.Range("E8").Value = "-6.15%"

# Row 9 (swapped with row 8 in the source data)
$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D9").Value = "0.9527"
$ws.Range("E9").Value = "0.52%"

# Row 10
$ws.Range("D10").Value = "0.1115"
$ws.Range("E10").Value = "-5.62%"

# Row 11
$ws.Range("D11").Value = "0.1860"
$ws.Range("E11").Value = "-1.79%"

# Row 12
$ws.Range("D12").Value = "10.75"
$ws.Range("E12").Value = "27.31%"

# Row 13
$ws.Range("D13").Value = "0.09806"
$ws.Range("E13").Value = "-0.76%"

# Row 14
$ws.Range("D14").Value = "0.04583"
$ws.Range("E14").Value = "9.69%"

# Row 15
$ws.Range("E15").Value = "0.17%"

# Row 16
$ws.Range("D16").Value = "0.001271"
$ws.Range("E16").Value = "-1.25%"

# Row 17
$ws.Range("E17").Value = "-4.08%"

# Row 18
$ws.Range("D18").Value = "0.005899"
$ws.Range("E18").Value = "-1.22%"

# Row 19
$ws.Range("D19").Value = "3.355"
$ws.Range("E19").Value = "-6.82%"

# Row 20
$ws.Range("D20").Value = "0.3474"
$ws.Range("E20").Value = "-0.35%"

# Row 21
$ws.Range("D21").Value = "0.1406"
$ws.Range("E21").Value = "2.29%"

# Row 22
$ws.Range("D22").Value = "0.2543"
$ws.Range("E22").Value = "0.19%"

# Row 23
$ws.Range("D23").Value = "0.001259"
$ws.Range("E23").Value = "1.39%"

# Row 24
$ws.Range("D24").Value = "0.004331"
$ws.Range("E24").Value = "-3.72%"

# Row 25
$ws.Range("D25").Value = "0.0001158"
$ws.Range("E25").Value = "-6.24%"

# Row 26
$ws.Range("D26").Value = "0.0003739"
$ws.Range("E26").Value = "-6.72%"

# Row 38
$ws.Range("D38").Value = "0.02559"
$ws.Range("E38").Value = "-2.83%"

# Row 39
$ws.Range("D39").Value = "0.05659"
$ws.Range("E39").Value = "3.11%"

# Row 40
$ws.Range("D40").Value = "0.007538"
$ws.Range("E40").Value = "-1.57%"

# Row 41
$ws.Range("D41").Value = "0.1395"
$ws.Range("E41").Value = "0.28%"

# Row 42
$ws.Range("D42").Value = "0.007617"
$ws.Range("E42").Value = "12.86%"

# Row 43
$ws.Range("D43").Value = "0.002013"
$ws.Range("E43").Value = "-2.09%"

# Row 44
$ws.Range("D44").Value = "0.008847"
$ws.Range("E44").Value = "-3.94%"

# Row 45
$ws.Range("D45").Value = "0.00007112"
$ws.Range("E45").Value = "-0.44%"

# Row 46
$ws.Range("E46").Value = "-0.59%"

# Row 47
$ws.Range("E47").Value = "54.72%"

# Row 48
$ws.Range("D48").Value = "0.003127"
$ws.Range("E48").Value = "-8.88%"

# Row 49
$ws.Range("D49").Value = "0.00002097"
$ws.Range("E49").Value = "-0.59%"

# Row 50
$ws.Range("D50").Value = "0.0001997"
$ws.Range("E50").Value = "-0.59%"
